{"js": "const replacements = [\n  [\"2025-02-18 Tuesday\", \"2025-02-19 Wednesday\"],\n  [\"59\u00f75=\", \"76\u00f73=\"],\n  [\"88\u00f72=\", \"89\u00f75=\"],\n  [\"52\u00f75=\", \"34\u00f72=\"],\n  [\"46\u00f77=\", \"11\u00f77=\"],\n  [\"29\u00f78=\", \"24\u00f76=\"],\n  [\"44\u00f79=\", \"95\u00f75=\"],\n  [\"99\u00f79=\", \"43\u00f75=\"],\n  [\"56\u00f78=\", \"90\u00f78=\"],\n  [\"53\u00f77=\", \"11\u00f72=\"],\n  [\"19\u00f76=\", \"40\u00f75=\"],\n  [\"96\u00f76=\", \"98\u00f77=\"],\n  [\"55\u00f75=\", \"12\u00f78=\"],\n  [\"82\u00f72=\", \"54\u00f78=\"],\n  [\"71\u00f73=\", \"16\u00f79=\"],\n  [\"13\u00f77=\", \"43\u00f74=\"],\n  [\"62\u00f74=\", \"73\u00f78=\"],\n  [\"36\u00f74=\", \"19\u00f77=\"],\n  [\"28\u00f73=\", \"43\u00f77=\"],\n  [\"96\u00f75=\", \"16\u00f74=\"],\n  [\"58\u00f74=\", \"90\u00f78=\"],\n  [\"48\u00f74=\", \"53\u00f77=\"],\n  [\"89\u00f77=\", \"29\u00f79=\"],\n  [\"54\u00f76=\", \"55\u00f76=\"],\n  [\"29\u00f72=\", \"29\u00f79=\"],\n  [\"33\u00f77=\", \"52\u00f75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-18 Tuesday\", \"2025-02-19 Wednesday\"),\n    @(\"59\u00f75=\", \"76\u00f73=\"),\n    @(\"88\u00f72=\", \"89\u00f75=\"),\n    @(\"52\u00f75=\", \"34\u00f72=\"),\n    @(\"46\u00f77=\", \"11\u00f77=\"),\n    @(\"29\u00f78=\", \"24\u00f76=\"),\n    @(\"44\u00f79=\", \"95\u00f75=\"),\n    @(\"99\u00f79=\", \"43\u00f75=\"),\n    @(\"56\u00f78=\", \"90\u00f78=\"),\n    @(\"53\u00f77=\", \"11\u00f72=\"),\n    @(\"19\u00f76=\", \"40\u00f75=\"),\n    @(\"96\u00f76=\", \"98\u00f77=\"),\n    @(\"55\u00f75=\", \"12\u00f78=\"),\n    @(\"82\u00f72=\", \"54\u00f78=\"),\n    @(\"71\u00f73=\", \"16\u00f79=\"),\n    @(\"13\u00f77=\", \"43\u00f74=\"),\n    @(\"62\u00f74=\", \"73\u00f78=\"),\n    @(\"36\u00f74=\", \"19\u00f77=\"),\n    @(\"28\u00f73=\", \"43\u00f77=\"),\n    @(\"96\u00f75=\", \"16\u00f74=\"),\n    @(\"58\u00f74=\", \"90\u00f78=\"),\n    @(\"48\u00f74=\", \"53\u00f77=\"),\n    @(\"89\u00f77=\", \"29\u00f79=\"),\n    @(\"54\u00f76=\", \"55\u00f76=\"),\n    @(\"29\u00f72=\", \"29\u00f79=\"),\n    @(\"33\u00f77=\", \"52\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}\n"}
